$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2023-12-02 -> 2023-12-03, i.e. 45262 -> 45263) for every data row
# (rows 2 through 27).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45263
}
